$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Stacking (SGD)"
$ws.Range("C2").Value = "Reduced"
$ws.Range("D2").Value = "Best"
$ws.Range("E2").Value = 0.7279608192341941
$ws.Range("F2").Value = 0.7448522829006267
$ws.Range("G2").Value = 0.7184801381692574
$ws.Range("H2").Value = 0.8129162559779494
$ws.Range("I2").Value = 832
$ws.Range("J2").Value = 803
$ws.Range("K2").Value = 326
$ws.Range("L2").Value = 285

# Row 3
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = "StackingCV"
$ws.Range("C3").Value = "Reduced"
$ws.Range("D3").Value = "Best"
$ws.Range("E3").Value = 0.7337488869100623
$ws.Range("F3").Value = 0.7403760071620412
$ws.Range("G3").Value = 0.7286343612334801
$ws.Range("H3").Value = 0.812694226357612
$ws.Range("I3").Value = 827
$ws.Range("J3").Value = 821
$ws.Range("K3").Value = 308
$ws.Range("L3").Value = 290

# Row 4
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "XGB"
$ws.Range("C4").Value = "Reduced"
$ws.Range("D4").Value = "Best"
$ws.Range("E4").Value = 0.7284060552092609
$ws.Range("F4").Value = 0.7600716204118174
$ws.Range("G4").Value = 0.7128463476070529
$ws.Range("H4").Value = 0.8112177293823691
$ws.Range("I4").Value = 849
$ws.Range("J4").Value = 787
$ws.Range("K4").Value = 342
$ws.Range("L4").Value = 268

# Row 5
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Voting"
$ws.Range("C5").Value = "Reduced"
$ws.Range("D5").Value = "Best"
$ws.Range("E5").Value = 0.7310774710596616
$ws.Range("F5").Value = 0.7287376902417189
$ws.Range("G5").Value = 0.7300448430493274
$ws.Range("H5").Value = 0.8066502629068595
$ws.Range("I5").Value = 814
$ws.Range("J5").Value = 828
$ws.Range("K5").Value = 301
$ws.Range("L5").Value = 303

# Row 6
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Stacking (SVC)"
$ws.Range("C6").Value = "Reduced"
$ws.Range("D6").Value = "Best"
$ws.Range("E6").Value = 0.7337488869100623
$ws.Range("F6").Value = 0.7394807520143241
$ws.Range("G6").Value = 0.7290379523389232
$ws.Range("H6").Value = 0.8044402752215738
$ws.Range("I6").Value = 826
$ws.Range("J6").Value = 822
$ws.Range("K6").Value = 307
$ws.Range("L6").Value = 291

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Random Forest"
$ws.Range("C7").Value = "Reduced"
$ws.Range("D7").Value = "Best"
$ws.Range("E7").Value = 0.719946571682992
$ws.Range("F7").Value = 0.7305282005371531
$ws.Range("G7").Value = 0.7132867132867133
$ws.Range("H7").Value = 0.7985997067623086
$ws.Range("I7").Value = 816
$ws.Range("J7").Value = 801
$ws.Range("K7").Value = 328
$ws.Range("L7").Value = 301

# Row 8
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Decision Tree"
$ws.Range("C8").Value = "Reduced"
$ws.Range("D8").Value = "Best"
$ws.Range("E8").Value = 0.695013357079252
$ws.Range("F8").Value = 0.6741271262309758
$ws.Range("G8").Value = 0.7011173184357542
$ws.Range("H8").Value = 0.7649689594661141
$ws.Range("I8").Value = 753
$ws.Range("J8").Value = 808
$ws.Range("K8").Value = 321
$ws.Range("L8").Value = 364

# Row 9
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "SVM (SVC)"
$ws.Range("C9").Value = "Full"
$ws.Range("D9").Value = "Default"
$ws.Range("E9").Value = 0.678539626001781
$ws.Range("F9").Value = 0.6454789615040286
$ws.Range("G9").Value = 0.6886341929321872
$ws.Range("H9").Value = 0.7510706981959301
$ws.Range("I9").Value = 721
$ws.Range("J9").Value = 803
$ws.Range("K9").Value = 326
$ws.Range("L9").Value = 396

# Row 10
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "Stacking (Logistic)"
$ws.Range("C10").Value = "Reduced"
$ws.Range("D10").Value = "Best"
$ws.Range("E10").Value = 0.649154051647373
$ws.Range("F10").Value = 0.5094001790510295
$ws.Range("G10").Value = 0.7033374536464772
$ws.Range("H10").Value = 0.740051288842298
$ws.Range("I10").Value = 569
$ws.Range("J10").Value = 889
$ws.Range("K10").Value = 240
$ws.Range("L10").Value = 548

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Vecstack"
$ws.Range("C11").Value = "Reduced"
$ws.Range("D11").Value = "Best"
$ws.Range("E11").Value = 0.7252894033837934
$ws.Range("F11").Value = 0.7600716204118174
$ws.Range("G11").Value = 0.7086811352253757
$ws.Range("H11").Value = 0.7353561553350941
$ws.Range("I11").Value = 849
$ws.Range("J11").Value = 780
$ws.Range("K11").Value = 349
$ws.Range("L11").Value = 268

# Row 12
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Logistic Regression"
$ws.Range("C12").Value = "Full"
$ws.Range("D12").Value = "Default"
$ws.Range("E12").Value = 0.6665182546749777
$ws.Range("F12").Value = 0.658012533572068
$ws.Range("G12").Value = 0.6669691470054446
$ws.Range("H12").Value = 0.7285434143239238
$ws.Range("I12").Value = 735
$ws.Range("J12").Value = 762
$ws.Range("K12").Value = 367
$ws.Range("L12").Value = 382

# Row 13
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Gaussian Naive-Bayes"
$ws.Range("C13").Value = "Full"
$ws.Range("D13").Value = "Default"
$ws.Range("E13").Value = 0.6567230632235085
$ws.Range("F13").Value = 0.5971351835273053
$ws.Range("G13").Value = 0.6751012145748988
$ws.Range("H13").Value = 0.7168908240708655
$ws.Range("I13").Value = 667
$ws.Range("J13").Value = 808
$ws.Range("K13").Value = 321
$ws.Range("L13").Value = 450

# Row 14
$ws.Range("A14").Value = 0
$ws.Range("B14").Value = "Bernoulli Naive-Bayes"
$ws.Range("C14").Value = "Full"
$ws.Range("D14").Value = "Default"
$ws.Range("E14").Value = 0.5894924309884239
$ws.Range("F14").Value = 0.3267681289167413
$ws.Range("G14").Value = 0.6822429906542056
$ws.Range("H14").Value = 0.6517199762428306
$ws.Range("I14").Value = 365
$ws.Range("J14").Value = 959
$ws.Range("K14").Value = 170
$ws.Range("L14").Value = 752
